# Apply scraped numeric updates to the per-class Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 9137.333000000001
$ws.Range("J48").Value = 9906.637000000001
$ws.Range("L48").Value = 29719.911
$ws.Range("N48").Value = -30303.911

$ws.Range("H56").Value = 9137.333000000001
$ws.Range("J56").Value = 9906.637000000001
$ws.Range("L56").Value = 29719.911
$ws.Range("N56").Value = -30787.911

$ws.Range("H106").Value = 8500.6
$ws.Range("I106").Value = 1519.9
$ws.Range("K106").Value = 1519.9
$ws.Range("M106").Value = -888.9000000000001

$ws.Range("H107").Value = 2450.2856
$ws.Range("I107").Value = 2358.8333
$ws.Range("K107").Value = 2358.8333
$ws.Range("M107").Value = -438.8332999999998

$ws.Range("H131").Value = 3251.7778
$ws.Range("I131").Value = 3614
$ws.Range("K131").Value = 10842
$ws.Range("M131").Value = -5802

$ws.Range("H136").Value = 69999
$ws.Range("J136").Value = 69999
$ws.Range("L136").Value = 69999
$ws.Range("N136").Value = -80199

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 14943.866
$ws.Range("I2").Value = 976.1429000000001
$ws.Range("J2").Value = 27165.625
$ws.Range("K2").Value = 976.1429000000001
$ws.Range("L2").Value = 27165.625
$ws.Range("M2").Value = -863.1429000000001
$ws.Range("N2").Value = -27391.625

$ws.Range("H13").Value = 6334.6665
$ws.Range("J13").Value = 7002
$ws.Range("L13").Value = 7002
$ws.Range("N13").Value = -7290

$ws.Range("H21").Value = 759.85
$ws.Range("I21").Value = 737.3570999999999
$ws.Range("J21").Value = 812.3333
$ws.Range("K21").Value = 737.3570999999999
$ws.Range("L21").Value = 812.3333
$ws.Range("M21").Value = -363.3570999999999
$ws.Range("N21").Value = -1560.3333

$ws.Range("H32").Value = 3069.5757
$ws.Range("I32").Value = 2642.3394
$ws.Range("K32").Value = 2642.3394
$ws.Range("M32").Value = -2355.3394

$ws.Range("H45").Value = 90913944
$ws.Range("I45").Value = 166668340
$ws.Range("K45").Value = 166668340
$ws.Range("M45").Value = -166667963

$ws.Range("H61").Value = 2764.6667
$ws.Range("I61").Value = 2764.6667
$ws.Range("K61").Value = 2764.6667
$ws.Range("M61").Value = -2552.6667

$ws.Range("H116").Value = 14943.866
$ws.Range("I116").Value = 976.1429000000001
$ws.Range("J116").Value = 27165.625
$ws.Range("K116").Value = 976.1429000000001
$ws.Range("L116").Value = 27165.625
$ws.Range("M116").Value = 1317.8571
$ws.Range("N116").Value = -31753.625

$ws.Range("H122").Value = 5332.6665
$ws.Range("I122").Value = 3838.8
$ws.Range("K122").Value = 11516.4
$ws.Range("M122").Value = -9066.400000000001

$ws.Range("H128").Value = 49597.8
$ws.Range("J128").Value = 49597.8
$ws.Range("L128").Value = 49597.8
$ws.Range("N128").Value = -59557.8

$ws.Range("H132").Value = 3065.6428
$ws.Range("I132").Value = 1838.2
$ws.Range("J132").Value = 3747.5557
$ws.Range("K132").Value = 5514.6
$ws.Range("L132").Value = 11242.6671
$ws.Range("M132").Value = -2984.6
$ws.Range("N132").Value = -16302.6671

$ws.Range("H136").Value = 2764.6667
$ws.Range("I136").Value = 2764.6667
$ws.Range("K136").Value = 8294.000100000001
$ws.Range("M136").Value = -5744.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 14943.866
$ws.Range("I3").Value = 976.1429000000001
$ws.Range("J3").Value = 27165.625
$ws.Range("K3").Value = 976.1429000000001
$ws.Range("L3").Value = 27165.625
$ws.Range("M3").Value = -862.1429000000001
$ws.Range("N3").Value = -27393.625

$ws.Range("H80").Value = 902.2143
$ws.Range("I80").Value = 903
$ws.Range("J80").Value = 901.9
$ws.Range("K80").Value = 903
$ws.Range("L80").Value = 901.9
$ws.Range("M80").Value = 95
$ws.Range("N80").Value = -2897.9

$ws.Range("H83").Value = 902.2143
$ws.Range("I83").Value = 903
$ws.Range("J83").Value = 901.9
$ws.Range("K83").Value = 4515
$ws.Range("L83").Value = 4509.5
$ws.Range("M83").Value = 477
$ws.Range("N83").Value = -14493.5

$ws.Range("H134").Value = 4185.722
$ws.Range("I134").Value = 2021.1333
$ws.Range("J134").Value = 15008.667
$ws.Range("K134").Value = 6063.3999
$ws.Range("L134").Value = 45026.001
$ws.Range("M134").Value = -3528.3999
$ws.Range("N134").Value = -50096.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 41759.57
$ws.Range("I31").Value = 2707.0527
$ws.Range("J31").Value = 124203.78
$ws.Range("K31").Value = 2707.0527
$ws.Range("L31").Value = 124203.78
$ws.Range("M31").Value = -2412.0527
$ws.Range("N31").Value = -124793.78

$ws.Range("H34").Value = 41759.57
$ws.Range("I34").Value = 2707.0527
$ws.Range("J34").Value = 124203.78
$ws.Range("K34").Value = 2707.0527
$ws.Range("L34").Value = 124203.78
$ws.Range("M34").Value = -2505.0527
$ws.Range("N34").Value = -124607.78

$ws.Range("H58").Value = 11986
$ws.Range("I58").Value = 5626
$ws.Range("J58").Value = 14711.714
$ws.Range("K58").Value = 5626
$ws.Range("L58").Value = 14711.714
$ws.Range("M58").Value = -5423
$ws.Range("N58").Value = -15117.714

$ws.Range("H107").Value = 966.6667
$ws.Range("I107").Value = 770.1818
$ws.Range("K107").Value = 770.1818
$ws.Range("M107").Value = 1149.8182

$ws.Range("H136").Value = 11986
$ws.Range("I136").Value = 5626
$ws.Range("J136").Value = 14711.714
$ws.Range("K136").Value = 16878
$ws.Range("L136").Value = 44135.142
$ws.Range("M136").Value = -14328
$ws.Range("N136").Value = -49235.142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 928.8
$ws.Range("I36").Value = 314
$ws.Range("J36").Value = 1851
$ws.Range("K36").Value = 942
$ws.Range("L36").Value = 5553
$ws.Range("M36").Value = -773
$ws.Range("N36").Value = -5891

$ws.Range("H40").Value = 414.2
$ws.Range("J40").Value = 512.625
$ws.Range("L40").Value = 2050.5
$ws.Range("N40").Value = -2188.5

$ws.Range("H63").Value = 12416.167
$ws.Range("I63").Value = 3796
$ws.Range("J63").Value = 18573.428
$ws.Range("K63").Value = 11388
$ws.Range("L63").Value = 55720.284
$ws.Range("M63").Value = -10639
$ws.Range("N63").Value = -57218.284

$ws.Range("H66").Value = 12416.167
$ws.Range("I66").Value = 3796
$ws.Range("J66").Value = 18573.428
$ws.Range("K66").Value = 34164
$ws.Range("L66").Value = 167160.852
$ws.Range("M66").Value = -30420
$ws.Range("N66").Value = -174648.852

$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("N78").ClearContents()

$ws.Range("H114").Value = 1289.8334
$ws.Range("I114").Value = 939.5714
$ws.Range("K114").Value = 2818.7142
$ws.Range("M114").Value = 435.2857999999997

$ws.Range("H131").Value = 27780580
$ws.Range("J131").Value = 32408824
$ws.Range("L131").Value = 97226472
$ws.Range("N131").Value = -97236552

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7106.2666
$ws.Range("I80").Value = 5208.9
$ws.Range("J80").Value = 10901
$ws.Range("K80").Value = 5208.9
$ws.Range("L80").Value = 10901
$ws.Range("M80").Value = -4210.9
$ws.Range("N80").Value = -12897

$ws.Range("H83").Value = 7106.2666
$ws.Range("I83").Value = 5208.9
$ws.Range("J83").Value = 10901
$ws.Range("K83").Value = 26044.5
$ws.Range("L83").Value = 54505
$ws.Range("M83").Value = -21052.5
$ws.Range("N83").Value = -64489

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H107").Value = 1892.9
$ws.Range("I107").Value = 776.4286
$ws.Range("K107").Value = 776.4286
$ws.Range("M107").Value = 1143.5714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 12098.444
$ws.Range("I7").Value = 2943.1667
$ws.Range("K7").Value = 2943.1667
$ws.Range("M7").Value = -2831.1667

$ws.Range("H22").Value = 7067.3335
$ws.Range("I22").Value = 1100
$ws.Range("K22").Value = 1100
$ws.Range("M22").Value = -805

$ws.Range("H27").Value = 7067.3335
$ws.Range("I27").Value = 1100
$ws.Range("K27").Value = 1100
$ws.Range("M27").Value = -993

$ws.Range("H46").Value = 4416.6665
$ws.Range("I46").Value = 3500
$ws.Range("J46").Value = 4875
$ws.Range("K46").Value = 3500
$ws.Range("L46").Value = 4875
$ws.Range("M46").Value = -3312
$ws.Range("N46").Value = -5251

$ws.Range("H126").Value = 12098.444
$ws.Range("I126").Value = 2943.1667
$ws.Range("K126").Value = 8829.500100000001
$ws.Range("M126").Value = -6359.500100000001

$ws.Range("H132").Value = 3670.3845
$ws.Range("I132").Value = 1903.9474
$ws.Range("K132").Value = 5711.8422
$ws.Range("M132").Value = -3181.8422

$ws.Range("H133").Value = 49263.2
$ws.Range("J133").Value = 49263.2
$ws.Range("L133").Value = 49263.2
$ws.Range("N133").Value = -54323.2

$ws.Range("H136").Value = 3750.88
$ws.Range("I136").Value = 1852.8182
$ws.Range("J136").Value = 17670
$ws.Range("K136").Value = 5558.4546
$ws.Range("L136").Value = 53010
$ws.Range("M136").Value = -3008.4546
$ws.Range("N136").Value = -58110

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 7462.9688
$ws.Range("I122").Value = 1874.0625
$ws.Range("K122").Value = 5622.1875
$ws.Range("M122").Value = -3172.1875

$ws.Range("H132").Value = 2725.6667
$ws.Range("I132").Value = 2572.3784
$ws.Range("K132").Value = 7717.135200000001
$ws.Range("M132").Value = -5187.135200000001

$ws.Range("H136").Value = 9244.647000000001
$ws.Range("I136").Value = 6258
$ws.Range("K136").Value = 18774
$ws.Range("M136").Value = -16224
